$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1424.25
$ws.Range("I18").Value = 634.2
$ws.Range("J18").Value = 1988.5714
$ws.Range("K18").Value = 634.2
$ws.Range("L18").Value = 1988.5714
$ws.Range("M18").Value = -350.2
$ws.Range("N18").Value = -2556.5714
$ws.Range("H28").Value = 477.9375
$ws.Range("I28").Value = 477.9375
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 477.9375
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = 7.0625
$ws.Range("H40").Value = 1612.3334
$ws.Range("I40").Value = 999.875
$ws.Range("J40").Value = 1918.5625
$ws.Range("K40").Value = 999.875
$ws.Range("L40").Value = 1918.5625
$ws.Range("M40").Value = -824.875
$ws.Range("N40").Value = -2268.5625
$ws.Range("H97").Value = 558.75
$ws.Range("J97").Value = 567.1429000000001
$ws.Range("L97").Value = 1701.4287
$ws.Range("N97").Value = -2693.4287
$ws.Range("H113").Value = 2940.4546
$ws.Range("I113").Value = 2810.8333
$ws.Range("J113").Value = 3096
$ws.Range("K113").Value = 2810.8333
$ws.Range("L113").Value = 3096
$ws.Range("M113").Value = 443.1667000000002
$ws.Range("N113").Value = -9604
$ws.Range("H132").Value = 13895963
$ws.Range("I132").Value = 7949.385
$ws.Range("J132").Value = 50004800
$ws.Range("K132").Value = 23848.155
$ws.Range("L132").Value = 150014400
$ws.Range("M132").Value = -21318.155
$ws.Range("N132").Value = -150019460
$ws.Range("H137").Value = 1109.5
$ws.Range("I137").Value = 1229
$ws.Range("J137").Value = 990
$ws.Range("K137").Value = 3687
$ws.Range("L137").Value = 2970
$ws.Range("M137").Value = -1137
$ws.Range("N137").Value = -8070
$ws.Range("H138").Value = 4446945.5
$ws.Range("J138").Value = 3568.0312
$ws.Range("L138").Value = 10704.0936
$ws.Range("N138").Value = -20984.0936
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 362.58334
$ws.Range("I5").Value = 331.375
$ws.Range("J5").Value = 425
$ws.Range("K5").Value = 331.375
$ws.Range("L5").Value = 425
$ws.Range("M5").Value = -219.375
$ws.Range("N5").Value = -649
$ws.Range("H6").Value = 28700
$ws.Range("I6").Value = 28700
$ws.Range("K6").Value = 28700
$ws.Range("M6").Value = -28527
$ws.Range("H32").Value = 8705.76
$ws.Range("I32").Value = 7215.747
$ws.Range("J32").Value = 15980.529
$ws.Range("K32").Value = 7215.747
$ws.Range("L32").Value = 15980.529
$ws.Range("M32").Value = -6928.747
$ws.Range("N32").Value = -16554.529
$ws.Range("H61").Value = 2822.1794
$ws.Range("I61").Value = 2945.6775
$ws.Range("J61").Value = 2343.625
$ws.Range("K61").Value = 2945.6775
$ws.Range("L61").Value = 2343.625
$ws.Range("M61").Value = -2733.6775
$ws.Range("N61").Value = -2767.625
$ws.Range("H74").Value = 2629.2
$ws.Range("I74").Value = 3256
$ws.Range("J74").Value = 2401.2727
$ws.Range("K74").Value = 3256
$ws.Range("L74").Value = 2401.2727
$ws.Range("M74").Value = -2382
$ws.Range("N74").Value = -4149.2727
$ws.Range("H77").Value = 2629.2
$ws.Range("I77").Value = 3256
$ws.Range("J77").Value = 2401.2727
$ws.Range("K77").Value = 16280
$ws.Range("L77").Value = 12006.3635
$ws.Range("M77").Value = -11912
$ws.Range("N77").Value = -20742.3635
$ws.Range("H107").Value = 28864
$ws.Range("J107").Value = 28864
$ws.Range("L107").Value = 28864
$ws.Range("N107").Value = -36544
$ws.Range("H110").Value = 1713.5
$ws.Range("I110").Value = 1238.7778
$ws.Range("K110").Value = 1238.7778
$ws.Range("M110").Value = 806.2221999999999
$ws.Range("H117").Value = 49646.4
$ws.Range("J117").Value = 49646.4
$ws.Range("L117").Value = 49646.4
$ws.Range("N117").Value = -58824.4
$ws.Range("H127").Value = 57477.273
$ws.Range("I127").Value = 57700
$ws.Range("J127").Value = 57455
$ws.Range("K127").Value = 57700
$ws.Range("L127").Value = 57455
$ws.Range("M127").Value = -52740
$ws.Range("N127").Value = -67375
$ws.Range("H132").Value = 5557538
$ws.Range("I132").Value = 7354531.5
$ws.Range("J132").Value = 3193.5454
$ws.Range("K132").Value = 22063594.5
$ws.Range("L132").Value = 9580.636200000001
$ws.Range("M132").Value = -22061064.5
$ws.Range("N132").Value = -14640.6362
$ws.Range("H133").Value = 32366
$ws.Range("J133").Value = 32366
$ws.Range("L133").Value = 32366
$ws.Range("N133").Value = -37426
$ws.Range("H136").Value = 2822.1794
$ws.Range("I136").Value = 2945.6775
$ws.Range("J136").Value = 2343.625
$ws.Range("K136").Value = 8837.032499999999
$ws.Range("L136").Value = 7030.875
$ws.Range("M136").Value = -6287.032499999999
$ws.Range("N136").Value = -12130.875
$ws.Range("H137").Value = 76922.22
$ws.Range("I137").Value = 39950
$ws.Range("J137").Value = 79880
$ws.Range("K137").Value = 39950
$ws.Range("L137").Value = 79880
$ws.Range("M137").Value = -34850
$ws.Range("N137").Value = -90080
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 362.58334
$ws.Range("I4").Value = 331.375
$ws.Range("J4").Value = 425
$ws.Range("K4").Value = 331.375
$ws.Range("L4").Value = 425
$ws.Range("M4").Value = -216.375
$ws.Range("N4").Value = -655
$ws.Range("H107").Value = 1343.6364
$ws.Range("I107").Value = 1209.9032
$ws.Range("J107").Value = 3416.5
$ws.Range("K107").Value = 1209.9032
$ws.Range("L107").Value = 3416.5
$ws.Range("M107").Value = 710.0968
$ws.Range("N107").Value = -7256.5
$ws.Range("H133").Value = 58571.43
$ws.Range("J133").Value = 58571.43
$ws.Range("L133").Value = 58571.43
$ws.Range("N133").Value = -68691.42999999999
$ws.Range("H134").Value = 4343.3486
$ws.Range("I134").Value = 3536
$ws.Range("J134").Value = 5045.391
$ws.Range("K134").Value = 10608
$ws.Range("L134").Value = 15136.173
$ws.Range("M134").Value = -8073
$ws.Range("N134").Value = -20206.173
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7111111
$ws.Range("I6").Value = 6750000
$ws.Range("K6").Value = 6750000
$ws.Range("M6").Value = -6749887
$ws.Range("H7").Value = 130
$ws.Range("I7").Value = 71
$ws.Range("J7").Value = 425
$ws.Range("K7").Value = 71
$ws.Range("L7").Value = 425
$ws.Range("M7").Value = 42
$ws.Range("N7").Value = -651
$ws.Range("H31").Value = 2859.5166
$ws.Range("I31").Value = 2246.121
$ws.Range("K31").Value = 2246.121
$ws.Range("M31").Value = -1951.121
$ws.Range("H34").Value = 2859.5166
$ws.Range("I34").Value = 2246.121
$ws.Range("K34").Value = 2246.121
$ws.Range("M34").Value = -2044.121
$ws.Range("H52").Value = 33393.08
$ws.Range("J52").Value = 33393.08
$ws.Range("L52").Value = 33393.08
$ws.Range("N52").Value = -33981.08
$ws.Range("H87").Value = 26529.1
$ws.Range("I87").Value = 15301
$ws.Range("J87").Value = 27776.666
$ws.Range("K87").Value = 15301
$ws.Range("L87").Value = 27776.666
$ws.Range("M87").Value = -14115
$ws.Range("N87").Value = -30148.666
$ws.Range("H90").Value = 26529.1
$ws.Range("I90").Value = 15301
$ws.Range("J90").Value = 27776.666
$ws.Range("K90").Value = 45903
$ws.Range("L90").Value = 83329.99800000001
$ws.Range("M90").Value = -39975
$ws.Range("N90").Value = -95185.99800000001
$ws.Range("H99").Value = 2109.65
$ws.Range("I99").Value = 1850
$ws.Range("K99").Value = 1850
$ws.Range("M99").Value = -352
$ws.Range("H126").Value = 2109.65
$ws.Range("I126").Value = 1850
$ws.Range("K126").Value = 5550
$ws.Range("M126").Value = -3080
$ws.Range("H127").Value = 53285.715
$ws.Range("J127").Value = 53285.715
$ws.Range("L127").Value = 53285.715
$ws.Range("N127").Value = -63205.715
$ws.Range("H134").Value = 645190.9
$ws.Range("I134").Value = 1930.7084
$ws.Range("K134").Value = 5792.1252
$ws.Range("M134").Value = -3257.1252
$ws.Range("H141").Value = 75660.5
$ws.Range("J141").Value = 80819.75999999999
$ws.Range("L141").Value = 80819.75999999999
$ws.Range("N141").Value = -91179.75999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 4200
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 4200
$ws.Range("K35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("M35").Value = 12600
$ws.Range("N35").Value = -13176
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1591.8182
$ws.Range("I107").Value = 2451
$ws.Range("J107").Value = 875.8333
$ws.Range("K107").Value = 2451
$ws.Range("L107").Value = 875.8333
$ws.Range("M107").Value = -531
$ws.Range("N107").Value = -4715.8333
$ws.Range("H122").Value = 3041.8948
$ws.Range("I122").Value = 2833.8
$ws.Range("J122").Value = 3273.111
$ws.Range("K122").Value = 8501.400000000001
$ws.Range("L122").Value = 9819.332999999999
$ws.Range("M122").Value = -6051.400000000001
$ws.Range("N122").Value = -14719.333
$ws.Range("H126").Value = 3457.5454
$ws.Range("I126").Value = 2563.111
$ws.Range("J126").Value = 4076.7693
$ws.Range("K126").Value = 7689.333
$ws.Range("L126").Value = 12230.3079
$ws.Range("M126").Value = -5219.333
$ws.Range("N126").Value = -17170.3079
$ws.Range("H132").Value = 5739.6875
$ws.Range("I132").Value = 6514.4
$ws.Range("J132").Value = 2972.8572
$ws.Range("K132").Value = 19543.2
$ws.Range("L132").Value = 8918.571599999999
$ws.Range("M132").Value = -17013.2
$ws.Range("N132").Value = -13978.5716
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4373.684
$ws.Range("I40").Value = 9698.333000000001
$ws.Range("K40").Value = 9698.333000000001
$ws.Range("M40").Value = -9562.333000000001
$ws.Range("H46").Value = 1383.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1383.5
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = 1383.5
$ws.Range("N46").Value = -1759.5
$ws.Range("H132").Value = 11117274
$ws.Range("I132").Value = 4252.613
$ws.Range("K132").Value = 12757.839
$ws.Range("M132").Value = -10227.839
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3048.3333
$ws.Range("I113").Value = 3463.3333
$ws.Range("J113").Value = 2633.3333
$ws.Range("K113").Value = 10389.9999
$ws.Range("L113").Value = 7899.999899999999
$ws.Range("M113").Value = -8219.999899999999
$ws.Range("N113").Value = -12239.9999
$ws.Range("H132").Value = 1986.8
$ws.Range("I132").Value = 1646.6757
$ws.Range("J132").Value = 2954.8462
$ws.Range("K132").Value = 4940.0271
$ws.Range("L132").Value = 8864.5386
$ws.Range("M132").Value = -2410.0271
$ws.Range("N132").Value = -13924.5386
